$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "Flair Fast_MCC scores"
$ws.Range("A12").Value = "Flair Fast_Classification time"
$ws.Range("A13").Value = "Flair Fast_Training time"

$ws.Range("A23").Value = "BertS_MCC scores"
$ws.Range("A24").Value = "BertS_Classification time"
$ws.Range("A25").Value = "BertS_Training time"

$ws.Range("A29").Value = "XLNetS_MCC scores"
$ws.Range("A30").Value = "XLNetS_Classification time"
$ws.Range("A31").Value = "XLNetS_Training time"

$ws.Range("A35").Value = "RobertaL_MCC scores"
$ws.Range("A36").Value = "RobertaL_Classification time"
$ws.Range("A37").Value = "RobertaL_Training time"

$ws.Range("A38").Value = "RobertaLS_MCC scores"
$ws.Range("A39").Value = "RobertaLS_Classification time"
$ws.Range("A40").Value = "RobertaLS_Training time"
